$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 186; this shifts existing rows 186..259 down to 187..260,
# matching the rest of the sheet (row 1 header + rows 2..185 stay untouched).
$ws.Rows.Item(186).Insert()

# Populate the newly inserted row 186 with the new data record.
$ws.Cells.Item(186, 1).Value = 9
$ws.Cells.Item(186, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(186, 3).Value = "Metropolitana"
$ws.Cells.Item(186, 4).Value = 44795
$ws.Cells.Item(186, 5).Value = 13
$ws.Cells.Item(186, 6).Value = 100112026
$ws.Cells.Item(186, 7).Value = "Haba"
$ws.Cells.Item(186, 8).Value = "Sin especificar"
$ws.Cells.Item(186, 9).Value = "Primera"
$ws.Cells.Item(186, 10).Value = 50
$ws.Cells.Item(186, 11).Value = 12000
$ws.Cells.Item(186, 12).Value = 12000
$ws.Cells.Item(186, 13).Value = 12000
$ws.Cells.Item(186, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(186, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(186, 16).Value = 480
$ws.Cells.Item(186, 17).Value = 25
$ws.Cells.Item(186, 18).Value = "Hortaliza"
